$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update K column (최종점수 / final score) values
$ws.Range("K2").Value = 68.7
$ws.Range("K3").Value = 68.7
$ws.Range("K4").Value = 64.1

# Update N column (MACRO_SCORE) values
$ws.Range("N2").Value = 85.82376350509293
$ws.Range("N3").Value = 85.82376350509293
$ws.Range("N4").Value = 85.82376350509293
